$d = $word.ActiveDocument

# Step 1: drop the standalone "Meta description" paragraph near the top.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# Step 2: find the closing AI-image-prompt paragraph. Insert a new bold
# "Play Bigger Bass Blizzard..." title paragraph right before it, then
# replace its own text with the review blurb (keeping the italic run).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Create an image featuring a happy Maya warrior*") {
        $p.Range.InsertParagraphBefore() | Out-Null

        # The blank paragraph just created sits at index $i; the original
        # (now "Create an image..." prompt) paragraph was pushed to $i + 1.
        $titlePara = $d.Paragraphs($i)
        $titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bigger Bass Blizzard - Christmas Catch Free | Review</w:t></w:r></w:p>"
        $titlePara.Range.InsertXML($titleXml) | Out-Null

        $descPara = $d.Paragraphs($i + 1)
        $descXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Bigger Bass Blizzard - Christmas Catch and play it for free. Enjoy stunning graphics, high maximum win, and a special bonus game.</w:t></w:r></w:p>"
        $descPara.Range.InsertXML($descXml) | Out-Null
        break
    }
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
